$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in the header (was 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = "5/24/2024"

# Update the prices for the three products listed
$ws.Range("D29").Value = 651
$ws.Range("D30").Value = 733
$ws.Range("D31").Value = 933
